$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.307.52"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.933.95"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.89"
$ws.Range("E5").Value = "  +4.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.64"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.730"
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000344"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.21"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "4.583.80"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.42"
$ws.Range("E14").Value = "  -3.06%  "
$ws.Range("D15").Value = "3.944.26"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("E18").Value = "  +7.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.80"
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").Value = "69.372.52"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.80"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.42"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.54"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.41"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.80"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.86"
$ws.Range("E26").Value = "  +6.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.10"
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.86"
$ws.Range("E28").Value = "  -5.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.65"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "707.29"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.32"
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.88"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "67.60"
$ws.Range("E34").Value = "  +11.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.444"
$ws.Range("E35").Value = "  +6.16%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.98"
$ws.Range("E37").Value = "  -6.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "40.59"
$ws.Range("E38").Value = "  -4.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.148"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0485"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.14"
$ws.Range("E43").Value = "  +7.72%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  -6.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.04"
$ws.Range("E45").Value = "  -7.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.143"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0359"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.37"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  +4.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("E51").Value = "  -1.90%  "
